# "wb: week 2 fixed on calories"
#
# Week 2 of the Meal Calendar sheet (rows 31-58) still had the old/broken
# Calories formula in column E (VLOOKUP against a dead #REF! range, keyed
# off the Meal column C). Week 1 (rows 2-29) had already been fixed to look
# up the Calories from the Recipes table using the Meal Name column (D)
# instead. Bring Week 2's column E formulas in line with Week 1, and
# restore the (previously missing) formula for E53 so the whole week sums
# correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meal Calendar")
$ws.Activate()

for ($row = 31; $row -le 58; $row++) {
    $cell = $ws.Cells.Item($row, 5)  # column E
    $cell.Formula = "=IFERROR(VLOOKUP(D$row, Recipes!`$B`$2:`$F`$1000, 2, FALSE), 0)"
}

# Restore the view state touched by the same save: scroll Week 2 into view
# and leave the selection on E58.
$ws.Range("E58").Select()
